$d = $word.ActiveDocument

function Remove-ParenBold($AnchorText, $BoldText) {
    # Locate the paragraph that contains the anchor text, and work with a
    # freshly constructed Document Range (Start/End) for every subsequent
    # Find so the search is strictly confined to that paragraph.
    $anchorRng = $d.Content.Duplicate
    $anchorRng.Find.Execute($AnchorText, $false, $false, $false, $false, $false, $true, 0, $false, "", 0)
    $paraRng = $anchorRng.Paragraphs(1).Range
    $pStart = $paraRng.Start
    $pEnd = $paraRng.End

    # Un-bold the bold text run first. Once it shares formatting with its
    # neighbouring (non-bold) runs, deleting the surrounding parentheses lets
    # the engine coalesce everything into a single run.
    $rngBold = $d.Range($pStart, $pEnd)
    $rngBold.Find.Execute($BoldText, $false, $false, $false, $false, $false, $true, 0, $false, "", 0)
    $rngBold.Font.Bold = 0

    # Delete the opening "(" immediately preceding the bold text.
    $rngOpen = $d.Range($pStart, $pEnd)
    $rngOpen.Find.Execute("(", $false, $false, $false, $false, $false, $true, 0, $false, "", 2)

    # Delete the closing ")" immediately following the bold text.
    $rngClose = $d.Range($pStart, $pEnd)
    $rngClose.Find.Execute(")", $false, $false, $false, $false, $false, $true, 0, $false, "", 2)
}

# --- Change 1: "Dates from 1987-01-01 to 2021-01-01 (Jan 2000 = 100)"
#     becomes  "Dates from 1987-01-01 to 2021-01-01 Jan 2000 = 100"
Remove-ParenBold "Dates from" "Jan 2000 = 100"

# --- Change 2: "Consumer Sentiment Index (Q1 1996 = 100)" becomes
#     "Consumer Sentiment Index Q1 1996 = 100"
Remove-ParenBold "Consumer Sentiment Index" "Q1 1996 = 100"

# The stray _GoBack bookmark that wrapped the "Q1 1996 = 100)" text is also
# removed as part of this revert.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
